$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 22: Cayton, 2005, Algorithms for manifold learning
$ws.Range("A22").Value = "cayton"
$ws.Range("B22").Value = 2005
$ws.Range("C22").Value = "algorithms for manifold learning"
$ws.Range("D22").Value = "paper"
$ws.Range("E22").Value = "review, meta analysis, dimensionality reduction"
$ws.Range("F22").Value = "definition of dimensionality reduction"

# New row 23: Verleysen, Francois, 2005, The curse of dimensionality in data mining and time series prediction
$ws.Range("A23").Value = "verleysen, francois"
$ws.Range("B23").Value = 2005
$ws.Range("C23").Value = "the curse of dimensionality in data mining and time series prediction"
$ws.Range("D23").Value = "paper"
$ws.Range("E23").Value = "cod"
$ws.Range("F23").Value = "intro"

# Copy formatting style from the row above (row 21) into the new rows
$ws.Range("A21:F21").Copy()
$ws.Range("A22:F23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Update the active cell selection to match the final state
$ws.Range("F23").Select() | Out-Null
